# Added Data for third preference (TP sheet), plus small update to SP sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the existing "SP" sheet: fill in the previously-blank M2 input
#    (3rd-pref "SA" percentage for the 19-response row) and move the
#    selection from M2 to M3. We do the selection update first while "SP"
#    is not yet the active sheet, so it won't end up marked as the
#    selected/active tab once we add & activate the new "TP" sheet below.
# ---------------------------------------------------------------------------
$sp = $wb.Worksheets.Item("SP")

$sp.Range("M3").Select()
$sp.Range("M2").Value = 0.1

# ---------------------------------------------------------------------------
# 2) Add the new "TP" (3rd preference) worksheet right after "SP".
# ---------------------------------------------------------------------------
$tp = $wb.Worksheets.Add($null, $sp)
$tp.Name = "TP"

# ---- Header row -------------------------------------------------------
$tp.Range("A1").Value = "No. of responses"
$tp.Range("B1").Value = "DIP"
$tp.Range("C1").Value = "AIS"
$tp.Range("D1").Value = "SEO"
$tp.Range("E1").Value = "SA"
$tp.Range("F1").Value = "UE"
$tp.Range("G1").Value = "ACN"
$tp.Range("H1").Value = "Sum"

$tp.Range("J1").Value = "DIP"
$tp.Range("K1").Value = "AIS"
$tp.Range("L1").Value = "SEO"
$tp.Range("M1").Value = "SA"
$tp.Range("N1").Value = "UE"
$tp.Range("O1").Value = "ACN"
$tp.Range("P1").Value = "Total"

# ---- Data rows (A = no. of responses; J:O = input percentages) --------
$tp.Range("A2").Value = 19
$tp.Range("J2").Value = 0
$tp.Range("K2").Value = 21.1
$tp.Range("L2").Value = 15.8
$tp.Range("M2").Value = 26.3
$tp.Range("N2").Value = 15.8
$tp.Range("O2").Value = 21.1

$tp.Range("A3").Value = 38
$tp.Range("J3").Value = 15.8
$tp.Range("K3").Value = 15.8
$tp.Range("L3").Value = 10.5
$tp.Range("M3").Value = 23.7
$tp.Range("N3").Value = 13.2
$tp.Range("O3").Value = 21.1

$tp.Range("A4").Value = 54
$tp.Range("J4").Value = 22.2
$tp.Range("K4").Value = 14.8
$tp.Range("L4").Value = 14.8
$tp.Range("M4").Value = 24.1
$tp.Range("N4").Value = 9.3
$tp.Range("O4").Value = 14.8

$tp.Range("A5").Value = 63
$tp.Range("J5").Value = 28.6
$tp.Range("K5").Value = 14.3
$tp.Range("L5").Value = 12.7
$tp.Range("M5").Value = 23.8
$tp.Range("N5").Value = 7.9
$tp.Range("O5").Value = 12.7

# ---- Formula columns B:G (percentage breakdown), one formula per cell -
foreach ($r in 2..5) {
    $tp.Range("B${r}").Formula = "=A${r}*J${r}%"
    $tp.Range("C${r}").Formula = "=A${r}*K${r}%"
    $tp.Range("D${r}").Formula = "=A${r}*L${r}%"
    $tp.Range("E${r}").Formula = "=A${r}*M${r}%"
    $tp.Range("F${r}").Formula = "=A${r}*N${r}%"
    $tp.Range("G${r}").Formula = "=A${r}*O${r}%"
}

# ---- Formula columns H and P (row sums) --------------------------------
# Row 2 is entered on its own, then rows 3:5 are set in one multi-cell
# assignment so they collapse into shared formulas on save -- matching how
# the rest of the workbook (e.g. the "SP" sheet) represents these sums.
$tp.Range("H2").Formula = "=SUM(B2:G2)"
$tp.Range("H3:H5").Formula = "=SUM(B3:G3)"

$tp.Range("P2").Formula = "=SUM(J2:O2)"
$tp.Range("P3:P5").Formula = "=SUM(J3:O3)"

# ---- Number formatting (integer display) for the computed columns -----
foreach ($r in 2..5) {
    $tp.Range("B${r}").NumberFormat = "0"
    $tp.Range("C${r}").NumberFormat = "0"
    $tp.Range("D${r}").NumberFormat = "0"
    $tp.Range("E${r}").NumberFormat = "0"
    $tp.Range("F${r}").NumberFormat = "0"
    $tp.Range("G${r}").NumberFormat = "0"
    $tp.Range("H${r}").NumberFormat = "0"
    $tp.Range("I${r}").NumberFormat = "0"
}

# ---- Selection for the new sheet (also becomes the active tab) --------
$tp.Range("O5").Select()
